$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (report volume/number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  40"
$ws.Range("C9").Value = "Report Covering the Week  10/2/2023  Through  10/8/2023"

# --- Crime statistics table updates (new weekly data) ---
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 2
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = -6.666666666666
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 14
$ws.Range("H16").Value = -7.142857142857
$ws.Range("I16").Value = 112
$ws.Range("J16").Value = 92
$ws.Range("K16").Value = 21.739130434782
$ws.Range("L16").Value = 41.772151898734
$ws.Range("M16").Value = -18.840579710144
$ws.Range("N16").Value = -80.984719864176
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 6
$ws.Range("E17").Value = -83.333333333333
$ws.Range("F17").Value = 16
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 168
$ws.Range("J17").Value = 136
$ws.Range("K17").Value = 23.529411764705
$ws.Range("L17").Value = 20.863309352518
$ws.Range("M17").Value = 78.723404255319
$ws.Range("N17").Value = -30
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = -33.333333333333
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = -38.095238095238
$ws.Range("I18").Value = 168
$ws.Range("J18").Value = 149
$ws.Range("K18").Value = 12.751677852349
$ws.Range("L18").Value = 54.128440366972
$ws.Range("M18").Value = -25
$ws.Range("N18").Value = -86.285714285714
$ws.Range("C19").Value = 10
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 77
$ws.Range("H19").Value = -41.558441558441
$ws.Range("I19").Value = 506
$ws.Range("J19").Value = 592
$ws.Range("K19").Value = -14.527027027027
$ws.Range("L19").Value = 15.261958997722
$ws.Range("M19").Value = 60.126582278481
$ws.Range("N19").Value = -14.814814814814
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = -50
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 16
$ws.Range("H20").Value = -43.75
$ws.Range("I20").Value = 122
$ws.Range("J20").Value = 115
$ws.Range("K20").Value = 6.086956521739
$ws.Range("L20").Value = 62.666666666666
$ws.Range("M20").Value = -4.6875
$ws.Range("N20").Value = -91.140159767610
$ws.Range("C21").Value = 17
$ws.Range("D21").Value = 32
$ws.Range("E21").Value = -46.875
$ws.Range("F21").Value = 98
$ws.Range("G21").Value = 146
$ws.Range("H21").Value = -32.876712328767
$ws.Range("I21").Value = 1092
$ws.Range("J21").Value = 1104
$ws.Range("K21").Value = -1.086956521739
$ws.Range("L21").Value = 27.272727272727
$ws.Range("M21").Value = 19.474835886214
$ws.Range("N21").Value = -73.037037037037
$ws.Range("C22").Value = "'0"
$ws.Range("L22").Value = 0
$ws.Range("C24").Value = 23
$ws.Range("D24").Value = 37
$ws.Range("E24").Value = -37.837837837837
$ws.Range("F24").Value = 99
$ws.Range("G24").Value = 171
$ws.Range("H24").Value = -42.105263157894
$ws.Range("I24").Value = 1359
$ws.Range("J24").Value = 1436
$ws.Range("K24").Value = -5.362116991643
$ws.Range("L24").Value = 46.601941747572
$ws.Range("M24").Value = 74.454428754813
$ws.Range("C25").Value = 19
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 58.333333333333
$ws.Range("F25").Value = 64
$ws.Range("G25").Value = 34
$ws.Range("H25").Value = 88.235294117647
$ws.Range("I25").Value = 468
$ws.Range("J25").Value = 349
$ws.Range("K25").Value = 34.097421203438
$ws.Range("L25").Value = 41.389728096676
$ws.Range("M25").Value = 37.647058823529
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -100
$ws.Range("F26").Value = 2
$ws.Range("G26").Value = 2
$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = -10
$ws.Range("C27").Value = "'0"
$ws.Range("E27").Value = -100
$ws.Range("J27").Value = 40
$ws.Range("K27").Value = 12.5
$ws.Range("L27").Value = 2.272727272727

# --- Fix up number formats/styles for cells that changed between N/A-text and numeric ---
$ws.Range("I14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("K14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("K14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("D26").PasteSpecial(-4122)
$ws.Range("K14").Copy()
$ws.Range("E26").PasteSpecial(-4122)
$ws.Range("I14").Copy()
$ws.Range("G26").PasteSpecial(-4122)
$ws.Range("K14").Copy()
$ws.Range("H26").PasteSpecial(-4122)
$ws.Range("A14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$excel.CutCopyMode = 0
